# Applies the weekly report "zeroing" update:
#  - Refresh the "Report Generated On" timestamp
#  - Zero out the Total Billed Amount summary cell
#  - Zero out every daily line-item Pricing value (column H) and the
#    per-day TOTAL rows, for each of the report's day sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Report Generated On" timestamp shown near the top of the report.
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"

# Zero out the Total Billed Amount in the summary block.
$ws.Range("C8").Value = 0

# Zero out the Pricing column (H) for every line item row and TOTAL row
# across all day sections (Thursday, Friday, Saturday blocks).
$zeroRows = @(16,17,18,19,20,21,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,55,56,57,58)

foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 8).Value = 0
}
